$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 16 (the review "I love this game because I feel
# that I am driving the car...") -- all rows below shift up by one.
$ws.Rows(16).Delete()

# The engine does not auto-shift hyperlink anchors when a row is deleted,
# so rebuild the hyperlinks collection to match the new row numbers.
$ws.Range("A1").Hyperlinks.Delete()

$links = @(
    @{Ref="D2"; Email="shmulmaor2@gmail.com"},
    @{Ref="C3"; Email="rocketaso@gmail.com"},
    @{Ref="D3"; Email="armonravid@gmail.com"},
    @{Ref="C5"; Email="ronoren61@gmail.com"},
    @{Ref="D5"; Email="nitanoren23@gmail.com"},
    @{Ref="C7"; Email="danfogel100@gmail.com"},
    @{Ref="D7"; Email="avishaybar12@gmail.com"},
    @{Ref="C8"; Email="danfogel100@gmail.com"},
    @{Ref="D8"; Email="avishaybar12@gmail.com"},
    @{Ref="D9"; Email="jorjkluni03@gmail.com"},
    @{Ref="C14"; Email="gazittalia1@gmail.com"},
    @{Ref="D14"; Email="hermanliran@gmail.com"},
    @{Ref="C16"; Email="budoyoni2@gmail.com"},
    @{Ref="D16"; Email="budoyoni@gmail.com"},
    @{Ref="C18"; Email="freelancernachus@gmail.com"},
    @{Ref="C21"; Email="itaisenior@gmail.com"},
    @{Ref="D21"; Email="vikicrestina@gmail.com"},
    @{Ref="C22"; Email="leviadlevi22@gmail.com"},
    @{Ref="D22"; Email="gazittalia1@gmail.com"},
    @{Ref="C23"; Email="nitanoren23@gmail.com"},
    @{Ref="D23"; Email="ronoren61@gmail.com"},
    @{Ref="C24"; Email="nachumshainc@gmail.com"},
    @{Ref="D24"; Email="nachushay@gmail.com"}
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Ref), "mailto:" + $link.Email, "", "", $link.Email)
}

# Restore the selection implied by the diff (it moves to A16 once the
# row holding the deleted review is gone).
$ws.Range("A16").Select() | Out-Null
